# Apply the "three order loop added" edit to the Ibkr-style workbook.
#
# Summary of the change:
#  - K226 changes from -1 to 0
#  - 23 new rows (227-249) are appended, each a clone of row 226's pattern
#    (same A/C/D/E/F/G/H/I/J values), where:
#      * rows 227-241 (15 rows) use Entry_Type "LIMIT"
#      * rows 242-249 (8 rows)  use Entry_Type "MARKET"
#      * all of the new rows have Activation (K) = 0, except the very last
#        row (249) which has Activation = 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, fix up the existing last row's Activation value.
$ws.Cells.Item(226, 11).Value = 0

# Clone row 226 (now with K=0) down through row 249, one row at a time so the
# formatting (number formats, styles, etc.) carries over just like a manual
# copy/paste of the row would do in Excel.
$srcRow = $ws.Range("A226:K226")
for ($r = 227; $r -le 249; $r++) {
    $destRow = $ws.Range("A$r`:K$r")
    $srcRow.Copy($destRow)
}

# Rows 227-241 keep Entry_Type "LIMIT" (inherited from the copy).
# Rows 242-249 switch Entry_Type to "MARKET".
$ws.Range("B242:B249").Value = "MARKET"

# The final new row (249) marks Activation = 1; all the other new rows stay 0.
$ws.Cells.Item(249, 11).Value = 1
